$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Câmera Canon)
$ws.Range("D2").Value = 5.2
$ws.Range("E2").Value = 5199.948
$ws.Range("G2").Value = 7279.9272

# Row 3 (Carro Renault)
$ws.Range("D3").Value = 5.65
$ws.Range("E3").Value = 25425
$ws.Range("G3").Value = 50850

# Row 4 (Notebook Dell)
$ws.Range("D4").Value = 5.2
$ws.Range("E4").Value = 4679.948
$ws.Range("G4").Value = 7955.9116

# Row 5 (IPhone)
$ws.Range("D5").Value = 5.2
$ws.Range("E5").Value = 4154.8
$ws.Range("G5").Value = 7063.16

# Row 6 (Carro Fiat)
$ws.Range("D6").Value = 5.65
$ws.Range("E6").Value = 16950
$ws.Range("G6").Value = 32205

# Row 7 (Celular Xiaomi)
$ws.Range("D7").Value = 5.2
$ws.Range("E7").Value = 2498.496
$ws.Range("G7").Value = 4996.992

# Row 8 (Joia 20g)
$ws.Range("D8").Value = 321.82
$ws.Range("E8").Value = 6436.4
$ws.Range("G8").Value = 7401.859999999999
